$wb = $excel.ActiveWorkbook

$wsBasePointers = $wb.Worksheets.Item("BasePointers")
$wsGameInfo = $wb.Worksheets.Item("GameInfo")

# BasePointers sheet: column F pointer offsets are stored as text
# (numeric-looking strings). Force text number format while writing so
# Excel doesn't auto-convert them to numeric cells, then restore the
# original (unstyled) look so no stray style gets left on the cell.

$rngF2 = $wsBasePointers.Range("F2")
$rngF2.NumberFormat = "@"
$rngF2.Value = "0"
$rngF2.Style = "Normal"

$rngF3 = $wsBasePointers.Range("F3")
$rngF3.NumberFormat = "@"
$rngF3.Value = "130958936"
$rngF3.Style = "Normal"

$rngF4 = $wsBasePointers.Range("F4")
$rngF4.NumberFormat = "@"
$rngF4.Value = "130959536"
$rngF4.Style = "Normal"

# GameInfo sheet: update the 2k26 build/version date string
$wsGameInfo.Range("F5").Value = "December 17th, 2025"
